$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Fermentation)
$ws.Range("B2").Value = 27.01457269758257
$ws.Range("D2").Value = 0.9343468331584022
$ws.Range("E2").Value = 0.363108018614811

# Row 3 (3-Phase Decanter)
$ws.Range("B3").Value = 3.398908506054629
$ws.Range("E3").Value = 0.6730299153427385

# Row 4 (Dehydration)
$ws.Range("B4").Value = 5.237513852805575
$ws.Range("C4").Value = 3.753292922668921
$ws.Range("D4").Value = 8.356435345475706

# Row 5 (Separation)
$ws.Range("B5").Value = 16.50675776278975
$ws.Range("C5").Value = 27.65653887058474
$ws.Range("D5").Value = 31.62661990543746
$ws.Range("E5").Value = 0.2799434290119038

# Row 6 (OSBL)
$ws.Range("B6").Value = 30.23887193470565
$ws.Range("C6").Value = 48.5621766103971
$ws.Range("E6").Value = 2.439800209241286
